# Update annotations for Ruilin
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B73 was stored as a text "3" - convert it to a true number 3
$ws.Range("B73").Value = 3

# Add new row 74 with the new annotation entry
$ws.Range("A74").Value = "Ruilin"

# B74 must remain a text value "2" (not a number) - force text formatting,
# assign it, then reset the cell style back to Normal so no extra
# number-format/style is left applied to the cell.
$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = "2"
$ws.Range("B74").Style = "Normal"

$ws.Range("C74").Value = "You do not understand the work by Veit et al."
$ws.Range("D74").Value = "CRT"
$ws.Range("E74").Value = "OTH"
$ws.Range("F74").Value = "9f35a425-2bea-4e69-9731-af889a0691d3"
$ws.Range("G74").Value = "r1Kr3TyAb_annotated.xlsx"
$ws.Range("H74").Value = "You do not understand the work by Veit et al."
